$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the trace-file reference in the explanatory note (row 10)
$ws.Range("A10").Value = "* Can be simulated using the 2-wide-opt.conf and tophat_vliw_scheduling_with_reordering.tr trace.  The result is going be for an in-order superscalar but the only difference for a VLIW is that instead of the hazard detection unit inserting bubbles, the compiler schedules instructions so that there are no hazards."

# The longer text needs a taller row to display comfortably
$ws.Rows.Item(10).RowHeight = 47.25

# Reflect the author's final selection/view state
$ws.Range("A10:Z10").Select() | Out-Null
